$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing row down into the new rows
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(344, 1).Value = 44418
$ws.Cells.Item(344, 2).Value = 0
$ws.Cells.Item(344, 3).Value = 11
$ws.Cells.Item(344, 4).Value = 62.61027946952018

$ws.Cells.Item(345, 1).Value = 44419
$ws.Cells.Item(345, 2).Value = 2
$ws.Cells.Item(345, 3).Value = 11
$ws.Cells.Item(345, 4).Value = 62.61027946952018

$ws.Cells.Item(346, 1).Value = 44420
$ws.Cells.Item(346, 2).Value = 1
$ws.Cells.Item(346, 3).Value = 12
$ws.Cells.Item(346, 4).Value = 68.30212305765838

$ws.Cells.Item(347, 1).Value = 44421
$ws.Cells.Item(347, 2).Value = 2
$ws.Cells.Item(347, 3).Value = 11
$ws.Cells.Item(347, 4).Value = 62.61027946952018

$ws.Cells.Item(348, 1).Value = 44422
$ws.Cells.Item(348, 2).Value = 6
$ws.Cells.Item(348, 3).Value = 13
$ws.Cells.Item(348, 4).Value = 73.99396664579658

$ws.Cells.Item(349, 1).Value = 44423
$ws.Cells.Item(349, 2).Value = 0
$ws.Cells.Item(349, 3).Value = 11
$ws.Cells.Item(349, 4).Value = 62.61027946952018

$ws.Cells.Item(350, 1).Value = 44424
$ws.Cells.Item(350, 2).Value = 3
$ws.Cells.Item(350, 3).Value = 14
$ws.Cells.Item(350, 4).Value = 79.68581023393477

$ws.Cells.Item(351, 1).Value = 44425
$ws.Cells.Item(351, 2).Value = 2
$ws.Cells.Item(351, 3).Value = 16
$ws.Cells.Item(351, 4).Value = 91.06949741021117

$ws.Cells.Item(352, 1).Value = 44426
$ws.Cells.Item(352, 2).Value = 1
$ws.Cells.Item(352, 3).Value = 15
$ws.Cells.Item(352, 4).Value = 85.37765382207297

$ws.Cells.Item(353, 1).Value = 44427
$ws.Cells.Item(353, 2).Value = 1
$ws.Cells.Item(353, 3).Value = 15
$ws.Cells.Item(353, 4).Value = 85.37765382207297

$ws.Cells.Item(354, 1).Value = 44428
$ws.Cells.Item(354, 2).Value = 4
$ws.Cells.Item(354, 3).Value = 17
$ws.Cells.Item(354, 4).Value = 96.76134099834937

$ws.Cells.Item(355, 1).Value = 44429
$ws.Cells.Item(355, 2).Value = 4
$ws.Cells.Item(355, 3).Value = 15
$ws.Cells.Item(355, 4).Value = 85.37765382207297

$ws.Cells.Item(356, 1).Value = 44430
$ws.Cells.Item(356, 2).Value = 1
$ws.Cells.Item(356, 3).Value = 16
$ws.Cells.Item(356, 4).Value = 91.06949741021117

$ws.Cells.Item(357, 1).Value = 44431
$ws.Cells.Item(357, 2).Value = 1
$ws.Cells.Item(357, 3).Value = 14
$ws.Cells.Item(357, 4).Value = 79.68581023393477
